$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain-text price strings (some look numeric, e.g. "534.89");
# prefix with an apostrophe so Excel stores them as text like the original file,
# then reset the style so no stray number-format / quote-prefix style is left behind.

$ws.Range("D2").Value = "'59.276.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("D3").Value = "'2.518.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.31%  '
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").Value = "'534.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = "'139.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.80%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  -1.72%  '
$ws.Range("D9").Value = "'2.524.15"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("E10").Value = '  +0.06%  '
$ws.Range("E11").Value = '  +1.18%  '
$ws.Range("D12").Value = "'5.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.91%  '
$ws.Range("D13").Value = "'0.356"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = "'2.968.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = "'23.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("D16").Value = "'59.196.99"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.47%  '
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = "'2.521.61"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.81%  '
$ws.Range("D19").Value = "'11.12"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.25%  '
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("D21").Value = "'325.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.29%  '
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").Value = "'5.82"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("D24").Value = "'63.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.36%  '
$ws.Range("E25").Value = '  -2.45%  '
$ws.Range("D26").Value = "'0.167"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("E28").Value = '  -2.16%  '
$ws.Range("D29").Value = "'6.92"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '
$ws.Range("D30").Value = "'0.0₃0777"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.11%  '
$ws.Range("E31").Value = '  -2.62%  '
$ws.Range("D32").Value = "'164.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.59%  '
$ws.Range("E33").Value = '  -2.04%  '
$ws.Range("E34").Value = '  +0.18%  '
$ws.Range("E35").Value = '  -9.98%  '
$ws.Range("D36").Value = "'18.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'4.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.05%  '
$ws.Range("E38").Value = '  -1.96%  '
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("E40").Value = '  -0.60%  '
$ws.Range("D41").Value = "'0.816"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.74%  '
$ws.Range("E42").Value = '  -7.34%  '
$ws.Range("D43").Value = "'279.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.87%  '
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("D46").Value = "'10.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.59%  '
$ws.Range("E47").Value = '  -0.15%  '
$ws.Range("D48").Value = "'123.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.18%  '
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("E50").Value = '  -1.70%  '
$ws.Range("D51").Value = "'17.75"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.42%  '
